$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resultats_merged")

$ws.Range("K3").Value = 169.84199999999998
$ws.Range("L3").Value = 208
$ws.Range("M3").Value = 17
$ws.Range("N3").Value = 60.252999999999929
$ws.Range("O3").Value = 5.4779999999999998
$ws.Range("Q3").Value = 0.20100000000002183
$ws.Range("R3").Value = 35.5
$ws.Range("S3").Value = 11
$ws.Range("K5").Value = 272.40999999999997
$ws.Range("L5").Value = 141
$ws.Range("M5").Value = 44
$ws.Range("R5").Value = 34
$ws.Range("K6").Value = 226.44399999999996
$ws.Range("L6").Value = 169
$ws.Range("M6").Value = 23
$ws.Range("N6").Value = 90.380000000000223
$ws.Range("O6").Value = 5.3159999999999998
$ws.Range("Q6").Value = 0.20299999999997453
$ws.Range("R6").Value = 39.9
$ws.Range("S6").Value = 17
$ws.Range("K7").Value = 495.83199999999999
$ws.Range("M7").Value = 36
$ws.Range("R7").Value = 23.5
$ws.Range("K8").Value = 288.45399999999995
$ws.Range("L8").Value = 160
$ws.Range("M8").Value = 39
$ws.Range("R8").Value = 30.4
$ws.Range("K9").Value = 226.40899999999999
$ws.Range("L9").Value = 188
$ws.Range("M9").Value = 26
$ws.Range("R9").Value = 23.8
$ws.Range("K10").Value = 61.499000000000024
$ws.Range("L10").Value = 7
$ws.Range("M10").Value = 13
$ws.Range("R10").Value = 52.5
$ws.Range("K11").Value = 280.44200000000001
$ws.Range("L11").Value = 132
$ws.Range("M11").Value = 36
$ws.Range("R11").Value = 38.700000000000003
$ws.Range("K12").Value = 219.04600000000005
$ws.Range("L12").Value = 120
$ws.Range("M12").Value = 30
$ws.Range("R12").Value = 38
$ws.Range("K13").Value = 446.64699999999993
$ws.Range("L13").Value = 118
$ws.Range("M13").Value = 52
$ws.Range("R13").Value = 44.3
$ws.Range("K17").Value = 291.48700000000002
$ws.Range("L17").Value = 106
$ws.Range("M17").Value = 36
$ws.Range("R17").Value = 31.4
$ws.Range("K18").Value = 236.37
$ws.Range("L18").Value = 120
$ws.Range("M18").Value = 19
$ws.Range("R18").Value = 35.700000000000003
$ws.Range("K19").Value = 429.322
$ws.Range("L19").Value = 163
$ws.Range("M19").Value = 57
$ws.Range("R19").Value = 34.9
$ws.Range("K20").Value = 230.46899999999999
$ws.Range("L20").Value = 127
$ws.Range("M20").Value = 26
$ws.Range("R20").Value = 55.2
$ws.Range("K21").Value = 234.67899999999997
$ws.Range("L21").Value = 125
$ws.Range("M21").Value = 20
$ws.Range("R21").Value = 52.6
$ws.Range("K23").Value = 231.47999999999996
$ws.Range("L23").Value = 131
$ws.Range("M23").Value = 33
$ws.Range("R23").Value = 50.2
$ws.Range("K25").Value = 198.25299999999993
$ws.Range("L25").Value = 103
$ws.Range("M25").Value = 24
$ws.Range("R25").Value = 52.3
$ws.Range("K26").Value = 201.68100000000004
$ws.Range("L26").Value = 176
$ws.Range("M26").Value = 25
$ws.Range("R26").Value = 35.200000000000003
$ws.Range("K28").Value = 260.60100000000011
$ws.Range("L28").Value = 164
$ws.Range("M28").Value = 24
$ws.Range("R28").Value = 24.1
$ws.Range("K29").Value = 222.04399999999998
$ws.Range("L29").Value = 165
$ws.Range("M29").Value = 39
$ws.Range("R29").Value = 40.299999999999997
$ws.Range("K30").Value = 230.92499999999995
$ws.Range("L30").Value = 106
$ws.Range("M30").Value = 26
$ws.Range("R30").Value = 24.2
$ws.Range("K32").Value = 1092.6420000000001
$ws.Range("L32").Value = 445
$ws.Range("M32").Value = 98
$ws.Range("R32").Value = 37.9
